# 20220621 - pins configuration changed.
# Update the MICRO pin labels in column F (rows 7-9 and 18-20) of the
# GLCD pin configuration sheet to reflect the new pinout, and move the
# active selection to H17 (matching the saved sheet view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-9 (D/I, R/W, E) used to be wired to the SPI pins (SPI_NSS,
# SPI_MOSI, SPI_SCK); they are now wired to direct GPIO pins.
$ws.Range("F7").Value = "PB12"
$ws.Range("F8").Value = "PB15"
$ws.Range("F9").Value = "PB13"

# Rows 18-19 (CS1, CS2) keep using PB8 / PB9 (unchanged values).
$ws.Range("F18").Value = "PB8"
$ws.Range("F19").Value = "PB9"

# Row 20 (RST) used to be wired to NRST; it now reuses PB15.
$ws.Range("F20").Value = "PB15"
$ws.Range("G20").Value = "RESET"

# Move the selection to H17, matching the sheet view saved with the change.
$ws.Range("H17").Select() | Out-Null
